{"js": "const replacements = [\n  [\"925\u00d72=1850\", \"826\u00d74=3304\"],\n  [\"890\u00d72=1780\", \"470\u00d77=3290\"],\n  [\"371\u00d72=742\", \"292\u00d75=1460\"],\n  [\"978\u00d77=6846\", \"657\u00d78=5256\"],\n  [\"471\u00d76=2826\", \"777\u00d77=5439\"],\n  [\"461\u00d73=1383\", \"253\u00d79=2277\"],\n  [\"112\u00d78=896\", \"731\u00d78=5848\"],\n  [\"837\u00d72=1674\", \"825\u00d74=3300\"],\n  [\"410\u00d78=3280\", \"186\u00d73=558\"],\n  [\"453\u00d78=3624\", \"631\u00d77=4417\"],\n  [\"154\u00d72=308\", \"341\u00d78=2728\"],\n  [\"480\u00d78=3840\", \"586\u00d75=2930\"],\n  [\"660\u00d78=5280\", \"181\u00d73=543\"],\n  [\"566\u00d74=2264\", \"209\u00d76=1254\"],\n  [\"639\u00d76=3834\", \"885\u00d74=3540\"],\n  [\"330\u00d72=660\", \"115\u00d72=230\"],\n  [\"233\u00d79=2097\", \"555\u00d75=2775\"],\n  [\"178\u00d76=1068\", \"824\u00d74=3296\"],\n  [\"701\u00d79=6309\", \"403\u00d76=2418\"],\n  [\"190\u00d76=1140\", \"651\u00d77=4557\"],\n  [\"262\u00d76=1572\", \"276\u00d78=2208\"],\n  [\"767\u00d76=4602\", \"772\u00d77=5404\"],\n  [\"819\u00d75=4095\", \"664\u00d79=5976\"],\n  [\"910\u00d79=8190\", \"323\u00d75=1615\"],\n  [\"347\u00d72=694\", \"929\u00d78=7432\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load('text');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"925\u00d72=1850\", \"826\u00d74=3304\"),\n    @(\"890\u00d72=1780\", \"470\u00d77=3290\"),\n    @(\"371\u00d72=742\", \"292\u00d75=1460\"),\n    @(\"978\u00d77=6846\", \"657\u00d78=5256\"),\n    @(\"471\u00d76=2826\", \"777\u00d77=5439\"),\n    @(\"461\u00d73=1383\", \"253\u00d79=2277\"),\n    @(\"112\u00d78=896\", \"731\u00d78=5848\"),\n    @(\"837\u00d72=1674\", \"825\u00d74=3300\"),\n    @(\"410\u00d78=3280\", \"186\u00d73=558\"),\n    @(\"453\u00d78=3624\", \"631\u00d77=4417\"),\n    @(\"154\u00d72=308\", \"341\u00d78=2728\"),\n    @(\"480\u00d78=3840\", \"586\u00d75=2930\"),\n    @(\"660\u00d78=5280\", \"181\u00d73=543\"),\n    @(\"566\u00d74=2264\", \"209\u00d76=1254\"),\n    @(\"639\u00d76=3834\", \"885\u00d74=3540\"),\n    @(\"330\u00d72=660\", \"115\u00d72=230\"),\n    @(\"233\u00d79=2097\", \"555\u00d75=2775\"),\n    @(\"178\u00d76=1068\", \"824\u00d74=3296\"),\n    @(\"701\u00d79=6309\", \"403\u00d76=2418\"),\n    @(\"190\u00d76=1140\", \"651\u00d77=4557\"),\n    @(\"262\u00d76=1572\", \"276\u00d78=2208\"),\n    @(\"767\u00d76=4602\", \"772\u00d77=5404\"),\n    @(\"819\u00d75=4095\", \"664\u00d79=5976\"),\n    @(\"910\u00d79=8190\", \"323\u00d75=1615\"),\n    @(\"347\u00d72=694\", \"929\u00d78=7432\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute(\n        $oldText,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    )\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}"}
